$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Add new "Problems" header column in A1 (new shared string) ---
$ws.Range("A1").Value2 = "Problems"
$ws.Range("A1").Font.Bold = $true

# --- 2. Un-bold the data cells in columns A and C (rows 2-4) ---
$ws.Range("A2").Font.Bold = $false
$ws.Range("A3").Font.Bold = $false
$ws.Range("A4").Font.Bold = $false
$ws.Range("C3").Font.Bold = $false
$ws.Range("C4").Font.Bold = $false

# --- 3. Row 4: add the second (iterative) Binary Search solution columns ---
# F4 keeps "O(log n)" (time) but switches from the "Good" style to "Neutral"
$ws.Range("F4").Value2 = "O(log n)"
$ws.Range("F4").Style = "Neutral"

# G4 is new: space for the iterative solution -> O(log n), "Good" style
$ws.Range("G4").Value2 = "O(log n)"
$ws.Range("G4").Style = "Good"

# H4 is new: holds what used to be in F4 -> O(1), "Good" style
$ws.Range("H4").Value2 = "O(1)"
$ws.Range("H4").Style = "Good"

# --- 4. Add cell comments explaining the two solutions ---
$cmt1 = $ws.Range("F4").AddComment()
$cmt1.Text("Gokul Chagalamarri Nippani:" + [char]10 + "If implemented recurresively because of stack.")

$cmt2 = $ws.Range("H4").AddComment()
$cmt2.Text("Gokul Chagalamarri Nippani:" + [char]10 + "If implemented iteratively")

# --- 5. Update the selected cell shown when the workbook is reopened ---
$ws.Range("D14").Select()
